# LULU.xlsx edit:
#   - Main!J2: unit price 229 -> 220 (dependent formulas J4, J7 recalc)
#   - Main sheet: active-cell selection moves from J3 to J2
#   - Model!S1:W1: extend the year header series (2031..2035) with
#     "=<prev col>1+1" formulas, matching the existing L1:R1 pattern
#   - Model!U13: discount/decay rate 0% -> -1%
# All other changed cells in the diff (row 10, row 21, U15/U16/U17, etc.)
# are downstream formula results that recalc automatically.

$wb = $excel.ActiveWorkbook

$wsMain  = $wb.Worksheets.Item("Main")
$wsModel = $wb.Worksheets.Item("Model")

# --- Main sheet -----------------------------------------------------
$wsMain.Range("J2").Value = 220

# Move the active selection to J2 (was J3) without changing which
# sheet tab is active (Model stays the active tab, as in the source).
[void]$wsMain.Range("J2").Select()
$wsModel.Activate()

# --- Model sheet ------------------------------------------------------
# Fill in the previously-empty year columns S1:W1 with the same
# "prior column + 1" formula used for L1:R1.
$wsModel.Range("S1").Formula = "=R1+1"
$wsModel.Range("T1").Formula = "=S1+1"
$wsModel.Range("U1").Formula = "=T1+1"
$wsModel.Range("V1").Formula = "=U1+1"
$wsModel.Range("W1").Formula = "=V1+1"

# Update the U13 driver assumption.
$wsModel.Range("U13").Value = -0.01
